$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells: Wins / Losses / Ties in AD1:AF1 ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the bold/centered/bordered header formatting already used by
# A1:AC1 (style index 1) by copying formats from an existing header cell
# instead of re-building the font/border/alignment from scratch (which
# would create a brand-new style entry rather than reusing the existing
# one).
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# --- Season record for every player row: Wins=59, Losses=102, Ties=0 ---
for ($row = 2; $row -le 52; $row++) {
    $ws.Cells.Item($row, 30).Value = 59   # AD - Wins
    $ws.Cells.Item($row, 31).Value = 102  # AE - Losses
    $ws.Cells.Item($row, 32).Value = 0    # AF - Ties
}
